$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update phone number and mellicode for row 2. A leading apostrophe forces
# Excel to keep these numeric-looking strings as text, preserving the
# leading zeros (same as a person typing '02938423984 into the cell).
$ws.Range("D2").Value = "'02938423984"
$ws.Range("E2").Value = "'0239482309"

# Row 3: clear phone number entirely, there is no longer a value there.
$ws.Range("D3").ClearContents()

# Row 3: mellicode becomes an empty text value (quote-prefixed blank entry, like typing ' then Enter)
$ws.Range("E3").Value = "'"
